# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N ("Outstanding"), shifting the existing N/O/P columns
# ("Outstanding" / heading / "Disbursement") one place to the right (to
# O/P/Q). Then make the "Repayment schedule" sheet the active sheet/tab,
# with cell K19 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column in front of column N (pushes N:P -> O:Q).
$ws.Columns("N:N").Insert() | Out-Null

# Make this sheet the active one, and select K19 on it (matches the
# workbook's saved view state after the edit).
$ws.Activate() | Out-Null
$ws.Range("K19").Select() | Out-Null
